$wb = $excel.ActiveWorkbook

# "mysprs" sheet (sheet1.xml) - first sheet in the workbook
$ws1 = $wb.Worksheets.Item(1)
# "mysprs2" sheet (sheet2.xml) - second sheet, stays the active tab
$ws2 = $wb.Worksheets.Item(2)

# Update the two changed cell values on the "mysprs" sheet:
#   B11: "Step1"  -> "Step1~"
#   B12: "Step2*" -> "Step2"
$ws1.Range("B11").Value = "Step1~"
$ws1.Range("B12").Value = "Step2"

# Move the selection on "mysprs" to B12 ...
$ws1.Range("B12").Select()
# ... but keep "mysprs2" as the active sheet/tab.
$ws2.Activate()
